# Applies the scheduled-runner price/profit refresh to each class sheet's
# Leve profit table (columns H:N) — current market prices changed, so the
# derived NQ/HQ price and profit columns are recomputed per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4083.3333
$ws.Range("I40").Value = 4530
$ws.Range("J40").Value = 1850
$ws.Range("K40").Value = 4530
$ws.Range("L40").Value = 1850
$ws.Range("M40").Value = -4355
$ws.Range("N40").Value = -2200
$ws.Range("H43").Value = 2033.375
$ws.Range("I43").Value = 896.7143
$ws.Range("J43").Value = 9990
$ws.Range("K43").Value = 896.7143
$ws.Range("L43").Value = 9990
$ws.Range("M43").Value = -827.7143
$ws.Range("N43").Value = -10128
$ws.Range("H98").Value = 4040.4443
$ws.Range("I98").Value = 3833.1538
$ws.Range("J98").Value = 4579.4
$ws.Range("K98").Value = 3833.1538
$ws.Range("L98").Value = 4579.4
$ws.Range("M98").Value = -2335.1538
$ws.Range("N98").Value = -7575.4
$ws.Range("H107").Value = 2176.3635
$ws.Range("I107").Value = 2384
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 2384
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = -464
$ws.Range("N107").Value = -3940
$ws.Range("H122").Value = 4040.4443
$ws.Range("I122").Value = 3833.1538
$ws.Range("J122").Value = 4579.4
$ws.Range("K122").Value = 11499.4614
$ws.Range("L122").Value = 13738.2
$ws.Range("M122").Value = -9049.4614
$ws.Range("N122").Value = -18638.2
$ws.Range("H132").Value = 223876.78
$ws.Range("I132").Value = 255756.11
$ws.Range("K132").Value = 767268.33
$ws.Range("M132").Value = -764738.33
$ws.Range("H141").Value = 2559.3
$ws.Range("I141").Value = 1949.75
$ws.Range("K141").Value = 5849.25
$ws.Range("M141").Value = -669.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4284633
$ws.Range("I32").Value = 10068.15
$ws.Range("K32").Value = 10068.15
$ws.Range("M32").Value = -9781.15
$ws.Range("H74").Value = 4887
$ws.Range("I74").Value = 5107.75
$ws.Range("J74").Value = 4298.3335
$ws.Range("K74").Value = 5107.75
$ws.Range("L74").Value = 4298.3335
$ws.Range("M74").Value = -4233.75
$ws.Range("N74").Value = -6046.3335
$ws.Range("H77").Value = 4887
$ws.Range("I77").Value = 5107.75
$ws.Range("J77").Value = 4298.3335
$ws.Range("K77").Value = 25538.75
$ws.Range("L77").Value = 21491.6675
$ws.Range("M77").Value = -21170.75
$ws.Range("N77").Value = -30227.6675
$ws.Range("H111").Value = 34999
$ws.Range("J111").Value = 34999
$ws.Range("L111").Value = 34999
$ws.Range("N111").Value = -43179
$ws.Range("H122").Value = 4205.923
$ws.Range("I122").Value = 3527.75
$ws.Range("K122").Value = 10583.25
$ws.Range("M122").Value = -8133.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14287246
$ws.Range("J20").Value = 1701
$ws.Range("L20").Value = 1701
$ws.Range("N20").Value = -2195
$ws.Range("H22").Value = 3614
$ws.Range("I22").Value = 338.5
$ws.Range("J22").Value = 10165
$ws.Range("K22").Value = 338.5
$ws.Range("L22").Value = 10165
$ws.Range("M22").Value = -165.5
$ws.Range("N22").Value = -10511
$ws.Range("H94").Value = 4445.6787
$ws.Range("I94").Value = 1861.9375
$ws.Range("J94").Value = 7890.6665
$ws.Range("K94").Value = 1861.9375
$ws.Range("L94").Value = 7890.6665
$ws.Range("M94").Value = -1410.9375
$ws.Range("N94").Value = -8792.666499999999
$ws.Range("H107").Value = 3579029.2
$ws.Range("I107").Value = 4352992
$ws.Range("K107").Value = 4352992
$ws.Range("M107").Value = -4351072
$ws.Range("H130").Value = 60600
$ws.Range("J130").Value = 60600
$ws.Range("L130").Value = 60600
$ws.Range("N130").Value = -70640
$ws.Range("H132").Value = 95696.75
$ws.Range("J132").Value = 95696.75
$ws.Range("L132").Value = 95696.75
$ws.Range("N132").Value = -105816.75
$ws.Range("H134").Value = 1258340.9
$ws.Range("I134").Value = 1702435.2
$ws.Range("K134").Value = 5107305.6
$ws.Range("M134").Value = -5104770.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 29418304
$ws.Range("I16").Value = 41671612
$ws.Range("K16").Value = 41671612
$ws.Range("M16").Value = -41671325
$ws.Range("H50").Value = 38727
$ws.Range("J50").Value = 38727
$ws.Range("L50").Value = 38727
$ws.Range("N50").Value = -39977
$ws.Range("H99").Value = 10587566
$ws.Range("I99").Value = 15879371
$ws.Range("J99").Value = 3957
$ws.Range("K99").Value = 15879371
$ws.Range("L99").Value = 3957
$ws.Range("M99").Value = -15877873
$ws.Range("N99").Value = -6953
$ws.Range("H107").Value = 1112.6666
$ws.Range("I107").Value = 1119.5
$ws.Range("J107").Value = 1099
$ws.Range("K107").Value = 1119.5
$ws.Range("L107").Value = 1099
$ws.Range("M107").Value = 800.5
$ws.Range("N107").Value = -4939
$ws.Range("H113").Value = 29418304
$ws.Range("I113").Value = 41671612
$ws.Range("K113").Value = 41671612
$ws.Range("M113").Value = -41669442
$ws.Range("H122").Value = 5988.5884
$ws.Range("I122").Value = 1907.0769
$ws.Range("K122").Value = 5721.2307
$ws.Range("M122").Value = -3271.2307
$ws.Range("H126").Value = 10587566
$ws.Range("I126").Value = 15879371
$ws.Range("J126").Value = 3957
$ws.Range("K126").Value = 47638113
$ws.Range("L126").Value = 11871
$ws.Range("M126").Value = -47635643
$ws.Range("N126").Value = -16811
$ws.Range("H132").Value = 10434.24
$ws.Range("I132").Value = 4848.0527
$ws.Range("J132").Value = 28123.834
$ws.Range("K132").Value = 14544.1581
$ws.Range("L132").Value = 84371.50199999999
$ws.Range("M132").Value = -12014.1581
$ws.Range("N132").Value = -89431.50199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1608.6
$ws.Range("I29").Value = 1962.75
$ws.Range("J29").Value = 192
$ws.Range("K29").Value = 5888.25
$ws.Range("L29").Value = 576
$ws.Range("M29").Value = -5611.25
$ws.Range("N29").Value = -1130
$ws.Range("H113").Value = 5695249.5
$ws.Range("I113").Value = 17083626
$ws.Range("J113").Value = 1061.5
$ws.Range("K113").Value = 51250878
$ws.Range("L113").Value = 3184.5
$ws.Range("M113").Value = -51248708
$ws.Range("N113").Value = -7524.5
$ws.Range("H121").Value = 15393.333
$ws.Range("J121").Value = 17344.277
$ws.Range("L121").Value = 52032.83099999999
$ws.Range("N121").Value = -54652.83099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 9563.532999999999
$ws.Range("J41").Value = 5006.875
$ws.Range("L41").Value = 5006.875
$ws.Range("N41").Value = -5716.875
$ws.Range("H70").Value = 8503
$ws.Range("J70").Value = 9999
$ws.Range("L70").Value = 9999
$ws.Range("N70").Value = -10539
$ws.Range("H73").Value = 8503
$ws.Range("J73").Value = 9999
$ws.Range("L73").Value = 9999
$ws.Range("N73").Value = -11871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9411.333000000001
$ws.Range("I7").Value = 5787.8
$ws.Range("J7").Value = 16658.4
$ws.Range("K7").Value = 5787.8
$ws.Range("L7").Value = 16658.4
$ws.Range("M7").Value = -5675.8
$ws.Range("N7").Value = -16882.4
$ws.Range("H40").Value = 4815.4116
$ws.Range("I40").Value = 3817.6924
$ws.Range("K40").Value = 3817.6924
$ws.Range("M40").Value = -3681.6924
$ws.Range("H46").Value = 29413720
$ws.Range("I46").Value = 1035.1428
$ws.Range("J46").Value = 50002600
$ws.Range("K46").Value = 1035.1428
$ws.Range("L46").Value = 50002600
$ws.Range("M46").Value = -847.1428000000001
$ws.Range("N46").Value = -50002976
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -25956
$ws.Range("H55").Value = 2925.9644
$ws.Range("I55").Value = 1390
$ws.Range("J55").Value = 7533.857
$ws.Range("K55").Value = 1390
$ws.Range("L55").Value = 7533.857
$ws.Range("M55").Value = -1217
$ws.Range("N55").Value = -7879.857
$ws.Range("H100").Value = 5017.625
$ws.Range("J100").Value = 2041
$ws.Range("L100").Value = 2041
$ws.Range("N100").Value = -3123
$ws.Range("H122").Value = 2224099.2
$ws.Range("I122").Value = 3332399
$ws.Range("K122").Value = 9997197
$ws.Range("M122").Value = -9994747
$ws.Range("H126").Value = 9411.333000000001
$ws.Range("I126").Value = 5787.8
$ws.Range("J126").Value = 16658.4
$ws.Range("K126").Value = 17363.4
$ws.Range("L126").Value = 49975.2
$ws.Range("M126").Value = -14893.4
$ws.Range("N126").Value = -54915.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 25000
$ws.Range("J5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("N5").Value = -25224
$ws.Range("H51").Value = 23866.857
$ws.Range("J51").Value = 24249.5
$ws.Range("L51").Value = 24249.5
$ws.Range("N51").Value = -25269.5
$ws.Range("H81").Value = 2014.2
$ws.Range("J81").Value = 4312.25
$ws.Range("L81").Value = 8624.5
$ws.Range("N81").Value = -10746.5
$ws.Range("H84").Value = 2014.2
$ws.Range("J84").Value = 4312.25
$ws.Range("L84").Value = 43122.5
$ws.Range("N84").Value = -53730.5
